$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix rows 389-391: column A was stored as text ("20082200" etc.) and
# should become a real number, matching the corrected "fedrollover" dates. ---
$ws.Range("A389").Value = 20082200
$ws.Range("B389").Value = 0
$ws.Range("C389").Value = 0

$ws.Range("A390").Value = 20082300
$ws.Range("B390").Value = 0
$ws.Range("C390").Value = 0

$ws.Range("A391").Value = 20082400
$ws.Range("B391").Value = 0
$ws.Range("C391").Value = 0

# --- New row 392 (numeric date) ---
$ws.Range("A392").Value = 20082500
$ws.Range("B392").Value = 120000000000
$ws.Range("C392").Value = 110000000000

# --- New rows 393-398: column A stored as text dates ---
$ws.Range("A393").Value = "'20082600"
$ws.Range("B393").Value = 0
$ws.Range("C393").Value = 0

$ws.Range("A394").Value = "'20082700"
$ws.Range("B394").Value = 165000000000
$ws.Range("C394").Value = 162000000000

$ws.Range("A395").Value = "'20082800"
$ws.Range("B395").Value = 22000000000
$ws.Range("C395").Value = 0

$ws.Range("A396").Value = "'20082900"
$ws.Range("B396").Value = 0
$ws.Range("C396").Value = 0

$ws.Range("A397").Value = "'20083000"
$ws.Range("B397").Value = 0
$ws.Range("C397").Value = 0

$ws.Range("A398").Value = "'20083100"
$ws.Range("B398").Value = 180000000000
$ws.Range("C398").Value = 100000000000
